$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "72.335.14"
$ws.Range("E2").Formula = "  -0.19%  "

$ws.Range("D3").Formula = "2.647.99"
$ws.Range("E3").Formula = "  +0.18%  "

$ws.Range("E4").Formula = "  +0.02%  "

$ws.Range("D5").Formula = "'591.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Formula = "  -1.99%  "

$ws.Range("D6").Formula = "'175.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Formula = "  -2.62%  "

$ws.Range("E7").Formula = "  +0.00%  "

$ws.Range("E8").Formula = "  -0.75%  "

$ws.Range("B9").Formula = "Dogecoin"
$ws.Range("C9").Formula = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Formula = "'0.173"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Formula = "  -2.16%  "

$ws.Range("B10").Formula = "LidoStakedEther"
$ws.Range("C10").Formula = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D10").Formula = "2.646.18"
$ws.Range("E10").Formula = "  +0.16%  "

$ws.Range("E11").Formula = "  +1.46%  "

$ws.Range("D12").Formula = "'0.357"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Formula = "  -0.33%  "

$ws.Range("E13").Formula = "  -1.61%  "

$ws.Range("D14").Formula = "3.133.04"
$ws.Range("E14").Formula = "  +0.21%  "

$ws.Range("E15").Formula = "  -2.14%  "

$ws.Range("D16").Formula = "72.215.32"
$ws.Range("E16").Formula = "  -0.19%  "

$ws.Range("D17").Formula = "'26.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Formula = "  -2.08%  "

$ws.Range("D18").Formula = "2.651.22"
$ws.Range("E18").Formula = "  +0.46%  "

$ws.Range("D19").Formula = "'12.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Formula = "  +1.98%  "

$ws.Range("D20").Formula = "'8.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Formula = "  +1.12%  "

$ws.Range("D21").Formula = "'370.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Formula = "  -2.55%  "

$ws.Range("D22").Formula = "'4.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Formula = "  -0.27%  "

$ws.Range("D23").Formula = "'2.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Formula = "  +0.73%  "

$ws.Range("D24").Formula = "'71.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Formula = "  -2.00%  "

$ws.Range("E25").Formula = "  +0.06%  "

$ws.Range("E26").Formula = "  -2.34%  "

$ws.Range("D27").Formula = "'9.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Formula = "  -3.09%  "

$ws.Range("D28").Formula = "2.780.60"
$ws.Range("E28").Formula = "  +0.03%  "

$ws.Range("E29").Formula = "  +0.49%  "

$ws.Range("D30").Formula = "0.0₃0963"
$ws.Range("E30").Formula = "  +0.40%  "

$ws.Range("D31").Formula = "'8.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Formula = "  -1.14%  "

$ws.Range("D32").Formula = "'502.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Formula = "  -4.30%  "

$ws.Range("E33").Formula = "  -1.66%  "

$ws.Range("E34").Formula = "  -0.66%  "

$ws.Range("D35").Formula = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Formula = "  +0.05%  "

$ws.Range("D36").Formula = "'161.71"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Formula = "  -1.90%  "

$ws.Range("E37").Formula = "  +3.80%  "

$ws.Range("E39").Formula = "  -1.02%  "

$ws.Range("E40").Formula = "  -2.64%  "

$ws.Range("E41").Formula = "  -0.04%  "

$ws.Range("E42").Formula = "  -6.45%  "

$ws.Range("D43").Formula = "'2.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Formula = "  -2.35%  "

$ws.Range("E44").Formula = "  -3.43%  "

$ws.Range("D45").Formula = "'0.329"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Formula = "  -1.31%  "

$ws.Range("D46").Formula = "'39.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Formula = "  -0.59%  "

$ws.Range("D47").Formula = "'153.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Formula = "  +1.68%  "

$ws.Range("E48").Formula = "  +1.65%  "

$ws.Range("D49").Formula = "'3.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Formula = "  -1.11%  "

$ws.Range("E50").Formula = "  -0.98%  "

$ws.Range("E51").Formula = "  -1.10%  "
